$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear rows that will be fully rewritten
$ws.Range("A2:AY2").ClearContents()
$ws.Range("A3:AY3").ClearContents()
$ws.Range("A4:AY4").ClearContents()
$ws.Range("A5:AY5").ClearContents()
$ws.Range("A6:AY6").ClearContents()
$ws.Range("A7:AY7").ClearContents()
$ws.Range("A8:AY8").ClearContents()
$ws.Range("A9:AY9").ClearContents()
$ws.Range("A10:AY10").ClearContents()
$ws.Range("A14:AY14").ClearContents()

# Row 2
$ws.Range("A2").Value = 100189556
$ws.Range("B2").Value = 89392
$ws.Range("C2").Value = 'Ovaliderad'
$ws.Range("D2").Value = 'NT'
$ws.Range("E2").Value = 1202
$ws.Range("F2").Value = 'Ullticka'
$ws.Range("G2").Value = 'Phellinidium ferrugineofuscum'
$ws.Range("H2").Value = '(P.Karst.) Fiasson & Niemelä'
$ws.Range("P2").Value = 'Hedemora, Dlr'
$ws.Range("Q2").Value = 565710.6320501693
$ws.Range("R2").Value = 6700394.867918631
$ws.Range("S2").Value = 25
$ws.Range("T2").Value = 'Dalarna'
$ws.Range("U2").Value = 'Hedemora'
$ws.Range("V2").Value = 'Dalarna'
$ws.Range("W2").Value = 'Husby'
$ws.Range("Y2").Value = '''2022-04-24'
$ws.Range("Z2").Value = '13:23'
$ws.Range("AA2").Value = '''2022-04-24'
$ws.Range("AB2").Value = '13:23'
$ws.Range("AD2").Value = $false
$ws.Range("AE2").Value = $false
$ws.Range("AG2").Value = $false
$ws.Range("AW2").Value = 'Philipp Weiss'
$ws.Range("AX2").Value = 'Philipp Weiss'

# Row 3
$ws.Range("A3").Value = 96334345
$ws.Range("B3").Value = 89356
$ws.Range("C3").Value = 'Ovaliderad'
$ws.Range("D3").Value = 'LC'
$ws.Range("E3").Value = 5447
$ws.Range("F3").Value = 'Vedticka'
$ws.Range("G3").Value = 'Fuscoporia viticola'
$ws.Range("H3").Value = '(Schwein.) Murrill'
$ws.Range("P3").Value = 'Hedemora, Dlr'
$ws.Range("Q3").Value = 565348.9174108147
$ws.Range("R3").Value = 6700212.542966858
$ws.Range("S3").Value = 25
$ws.Range("T3").Value = 'Dalarna'
$ws.Range("U3").Value = 'Hedemora'
$ws.Range("V3").Value = 'Dalarna'
$ws.Range("W3").Value = 'Husby'
$ws.Range("Y3").Value = '''2021-09-27'
$ws.Range("Z3").Value = '11:22'
$ws.Range("AA3").Value = '''2021-09-27'
$ws.Range("AB3").Value = '11:22'
$ws.Range("AC3").Value = 'På granlåga i granplantage.'
$ws.Range("AD3").Value = $false
$ws.Range("AE3").Value = $false
$ws.Range("AG3").Value = $false
$ws.Range("AW3").Value = 'Annelie Hilmerby'
$ws.Range("AX3").Value = 'Annelie Hilmerby'

# Row 4
$ws.Range("A4").Value = 112425213
$ws.Range("B4").Value = 90466
$ws.Range("C4").Value = 'Ovaliderad'
$ws.Range("D4").Value = 'LC'
$ws.Range("E4").Value = 4769
$ws.Range("F4").Value = 'Svavelriska'
$ws.Range("G4").Value = 'Lactarius scrobiculatus'
$ws.Range("H4").Value = '(Scop.:Fr.) Fr.'
$ws.Range("P4").Value = 'Hedemora (Hedemora), Dlr'
$ws.Range("Q4").Value = 565673
$ws.Range("R4").Value = 6700080
$ws.Range("S4").Value = 15
$ws.Range("T4").Value = 'Dalarna'
$ws.Range("U4").Value = 'Hedemora'
$ws.Range("V4").Value = 'Dalarna'
$ws.Range("W4").Value = 'Husby'
$ws.Range("Y4").Value = '''2023-09-30'
$ws.Range("Z4").Value = '17:50'
$ws.Range("AA4").Value = '''2023-09-30'
$ws.Range("AB4").Value = '17:50'
$ws.Range("AD4").Value = $false
$ws.Range("AE4").Value = $false
$ws.Range("AG4").Value = $false
$ws.Range("AW4").Value = 'Philipp Weiss'
$ws.Range("AX4").Value = 'Philipp Weiss'

# Row 5
$ws.Range("A5").Value = 94360483
$ws.Range("B5").Value = 96334
$ws.Range("C5").Value = 'Ovaliderad'
$ws.Range("D5").Value = 'VU'
$ws.Range("E5").Value = 220787
$ws.Range("F5").Value = 'Knärot'
$ws.Range("G5").Value = 'Goodyera repens'
$ws.Range("H5").Value = '(L.) R. Br.'
$ws.Range("K5").Value = 'blomknopp'
$ws.Range("P5").Value = 'Sörboberget, Dlr'
$ws.Range("Q5").Value = 565784.6830013145
$ws.Range("R5").Value = 6699802.648075164
$ws.Range("S5").Value = 25
$ws.Range("T5").Value = 'Dalarna'
$ws.Range("U5").Value = 'Hedemora'
$ws.Range("V5").Value = 'Dalarna'
$ws.Range("W5").Value = 'Husby'
$ws.Range("Y5").Value = '''2021-06-20'
$ws.Range("Z5").Value = '00:00'
$ws.Range("AA5").Value = '''2021-06-20'
$ws.Range("AB5").Value = '00:00'
$ws.Range("AD5").Value = $false
$ws.Range("AE5").Value = $false
$ws.Range("AG5").Value = $false
$ws.Range("AW5").Value = 'Philipp Weiss'
$ws.Range("AX5").Value = 'Philipp Weiss'

# Row 6
$ws.Range("A6").Value = 94360735
$ws.Range("B6").Value = 98520
$ws.Range("C6").Value = 'Ovaliderad'
$ws.Range("D6").Value = 'LC'
$ws.Range("E6").Value = 222498
$ws.Range("F6").Value = 'Blåsippa'
$ws.Range("G6").Value = 'Hepatica nobilis'
$ws.Range("H6").Value = 'Schreb.'
$ws.Range("P6").Value = 'Sörboberget, Dlr'
$ws.Range("Q6").Value = 565681.6469905056
$ws.Range("R6").Value = 6699892.635672216
$ws.Range("S6").Value = 25
$ws.Range("T6").Value = 'Dalarna'
$ws.Range("U6").Value = 'Hedemora'
$ws.Range("V6").Value = 'Dalarna'
$ws.Range("W6").Value = 'Husby'
$ws.Range("Y6").Value = '''2021-06-20'
$ws.Range("Z6").Value = '00:00'
$ws.Range("AA6").Value = '''2021-06-20'
$ws.Range("AB6").Value = '00:00'
$ws.Range("AD6").Value = $false
$ws.Range("AE6").Value = $false
$ws.Range("AG6").Value = $false
$ws.Range("AH6").Value = 'Barrträd'
$ws.Range("AW6").Value = 'Philipp Weiss'
$ws.Range("AX6").Value = 'Philipp Weiss'

# Row 7
$ws.Range("A7").Value = 96241769
$ws.Range("B7").Value = 96334
$ws.Range("C7").Value = 'Ovaliderad'
$ws.Range("D7").Value = 'VU'
$ws.Range("E7").Value = 220787
$ws.Range("F7").Value = 'Knärot'
$ws.Range("G7").Value = 'Goodyera repens'
$ws.Range("H7").Value = '(L.) R. Br.'
$ws.Range("I7").Value = '''12'
$ws.Range("K7").Value = 'fullt utvecklade blad'
$ws.Range("P7").Value = 'Hedemora, Dlr'
$ws.Range("Q7").Value = 565750.3740434679
$ws.Range("R7").Value = 6699921.536218314
$ws.Range("S7").Value = 25
$ws.Range("T7").Value = 'Dalarna'
$ws.Range("U7").Value = 'Hedemora'
$ws.Range("V7").Value = 'Dalarna'
$ws.Range("W7").Value = 'Husby'
$ws.Range("Y7").Value = '''2021-09-22'
$ws.Range("Z7").Value = '16:26'
$ws.Range("AA7").Value = '''2021-09-22'
$ws.Range("AB7").Value = '16:26'
$ws.Range("AD7").Value = $false
$ws.Range("AE7").Value = $false
$ws.Range("AG7").Value = $false
$ws.Range("AW7").Value = 'Philipp Weiss'
$ws.Range("AX7").Value = 'Philipp Weiss'

# Row 8
$ws.Range("A8").Value = 96241731
$ws.Range("B8").Value = 77506
$ws.Range("C8").Value = 'Ovaliderad'
$ws.Range("D8").Value = 'NT'
$ws.Range("E8").Value = 6425
$ws.Range("F8").Value = 'Garnlav'
$ws.Range("G8").Value = 'Alectoria sarmentosa'
$ws.Range("H8").Value = '(Ach.) Ach.'
$ws.Range("P8").Value = 'Hedemora, Dlr'
$ws.Range("Q8").Value = 565750.3740434679
$ws.Range("R8").Value = 6699921.536218314
$ws.Range("S8").Value = 25
$ws.Range("T8").Value = 'Dalarna'
$ws.Range("U8").Value = 'Hedemora'
$ws.Range("V8").Value = 'Dalarna'
$ws.Range("W8").Value = 'Husby'
$ws.Range("Y8").Value = '''2021-09-22'
$ws.Range("Z8").Value = '16:25'
$ws.Range("AA8").Value = '''2021-09-22'
$ws.Range("AB8").Value = '16:25'
$ws.Range("AD8").Value = $false
$ws.Range("AE8").Value = $false
$ws.Range("AG8").Value = $false
$ws.Range("AW8").Value = 'Philipp Weiss'
$ws.Range("AX8").Value = 'Philipp Weiss'

# Row 9
$ws.Range("A9").Value = 96334977
$ws.Range("B9").Value = 93044
$ws.Range("C9").Value = 'Ovaliderad'
$ws.Range("D9").Value = 'LC'
$ws.Range("E9").Value = 2809
$ws.Range("F9").Value = 'Mörk husmossa'
$ws.Range("G9").Value = 'Hylocomiastrum umbratum'
$ws.Range("H9").Value = '(Hedw.) M.Fleisch.'
$ws.Range("P9").Value = 'Hedemora, Dlr'
$ws.Range("Q9").Value = 565546.2447679342
$ws.Range("R9").Value = 6699940.062954916
$ws.Range("S9").Value = 25
$ws.Range("T9").Value = 'Dalarna'
$ws.Range("U9").Value = 'Hedemora'
$ws.Range("V9").Value = 'Dalarna'
$ws.Range("W9").Value = 'Husby'
$ws.Range("Y9").Value = '''2021-09-27'
$ws.Range("Z9").Value = '11:50'
$ws.Range("AA9").Value = '''2021-09-27'
$ws.Range("AB9").Value = '11:50'
$ws.Range("AD9").Value = $false
$ws.Range("AE9").Value = $false
$ws.Range("AG9").Value = $false
$ws.Range("AW9").Value = 'Annelie Hilmerby'
$ws.Range("AX9").Value = 'Annelie Hilmerby'

# Row 10
$ws.Range("A10").Value = 96335108
$ws.Range("B10").Value = 90319
$ws.Range("C10").Value = 'Ovaliderad'
$ws.Range("D10").Value = 'LC'
$ws.Range("E10").Value = 4769
$ws.Range("F10").Value = 'Svavelriska'
$ws.Range("G10").Value = 'Lactarius scrobiculatus'
$ws.Range("H10").Value = '(Scop.:Fr.) Fr.'
$ws.Range("I10").Value = '''1'
$ws.Range("J10").Value = 'mycel'
$ws.Range("P10").Value = 'Hedemora, Dlr'
$ws.Range("Q10").Value = 565550.8325909179
$ws.Range("R10").Value = 6699959.899156544
$ws.Range("S10").Value = 25
$ws.Range("T10").Value = 'Dalarna'
$ws.Range("U10").Value = 'Hedemora'
$ws.Range("V10").Value = 'Dalarna'
$ws.Range("W10").Value = 'Husby'
$ws.Range("Y10").Value = '''2021-09-27'
$ws.Range("Z10").Value = '11:57'
$ws.Range("AA10").Value = '''2021-09-27'
$ws.Range("AB10").Value = '11:57'
$ws.Range("AD10").Value = $false
$ws.Range("AE10").Value = $false
$ws.Range("AG10").Value = $false
$ws.Range("AW10").Value = 'Annelie Hilmerby'
$ws.Range("AX10").Value = 'Annelie Hilmerby'

# Row 14
$ws.Range("A14").Value = 94360456
$ws.Range("B14").Value = 96334
$ws.Range("C14").Value = 'Ovaliderad'
$ws.Range("D14").Value = 'VU'
$ws.Range("E14").Value = 220787
$ws.Range("F14").Value = 'Knärot'
$ws.Range("G14").Value = 'Goodyera repens'
$ws.Range("H14").Value = '(L.) R. Br.'
$ws.Range("P14").Value = 'Hedemora, Dlr'
$ws.Range("Q14").Value = 565789.2480928447
$ws.Range("R14").Value = 6699769.148744237
$ws.Range("S14").Value = 25
$ws.Range("T14").Value = 'Dalarna'
$ws.Range("U14").Value = 'Hedemora'
$ws.Range("V14").Value = 'Dalarna'
$ws.Range("W14").Value = 'Husby'
$ws.Range("Y14").Value = '''2021-06-20'
$ws.Range("Z14").Value = '09:57'
$ws.Range("AA14").Value = '''2021-06-20'
$ws.Range("AB14").Value = '09:57'
$ws.Range("AD14").Value = $false
$ws.Range("AE14").Value = $false
$ws.Range("AG14").Value = $false
$ws.Range("AW14").Value = 'Annelie Hilmerby'
$ws.Range("AX14").Value = 'Annelie Hilmerby'

# Row 12: update B only
$ws.Range("B12").Value = 98961
# Row 13: update B only
$ws.Range("B13").Value = 89539

# Normalize already-empty cells in untouched rows to stay fully empty
$ws.Range("I11").ClearContents()
$ws.Range("AT11").ClearContents()
$ws.Range("AY11").ClearContents()
$ws.Range("AT12").ClearContents()
$ws.Range("AY12").ClearContents()
$ws.Range("AT13").ClearContents()
$ws.Range("AY13").ClearContents()

Write-Host "done"